$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$oldText = $wsHoja1.Range("A1").Value2
$newText = $oldText -replace "1000 Bs = 2\.35 = 8880\.75 pesos", "1000 Bs = 2.33 = 8819.69 pesos"
$newText = $newText -replace "8880\.75 pesos = 2\.33 = 946\.35 Bs", "8819.69 pesos = 2.32 = 936.13 Bs"
$wsHoja1.Range("A1").Value2 = $newText

# --- Sheet "tasas": update the rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 428.7
$wsTasas.Range("O10").Value = 3781
$wsTasas.Range("N12").Value = 3808.99
$wsTasas.Range("O12").Value = 404.288
